# Atualizado por script em 05-11-2023 08:45
# Adds the two new match rows (26 and 27) to the Gibraltar National League
# sheet, matching the style (bold/bordered index column, date-formatted
# kickoff column) already used by the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (25) down into the two
# new rows so the new cells inherit the same styles (s="1" on column A,
# s="2" on column E) as the rest of the table.
$ws.Range("A25:V25").Copy()
$ws.Range("A26:V26").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A27:V27").PasteSpecial(-4122) # xlPasteFormats

function Set-MatchRow(
    $Row,
    $Indice,
    $Pais,
    $Torneio,
    $Temporada,
    $DataPartida,
    $Home,
    $HomeGols,
    $Away,
    $AwayGols,
    $HomeOpeningOdds,
    $HomeOpeningDataHora,
    $HomeClosingOdds,
    $HomeClosingDataHora,
    $DrawOpeningOdds,
    $DrawOpeningDataHora,
    $DrawClosingOdds,
    $DrawClosingDataHora,
    $AwayOpeningOdds,
    $AwayOpeningDataHora,
    $AwayClosingOdds,
    $AwayClosingDataHora,
    $Url
) {
    $ws.Cells.Item($Row, 1).Value = $Indice
    $ws.Cells.Item($Row, 2).Value = $Pais
    $ws.Cells.Item($Row, 3).Value = $Torneio
    $ws.Cells.Item($Row, 4).Value = $Temporada
    $ws.Cells.Item($Row, 5).Value = $DataPartida
    $ws.Cells.Item($Row, 6).Value = $Home
    $ws.Cells.Item($Row, 7).Value = $HomeGols
    $ws.Cells.Item($Row, 8).Value = $Away
    $ws.Cells.Item($Row, 9).Value = $AwayGols
    $ws.Cells.Item($Row, 10).Value = $HomeOpeningOdds
    $ws.Cells.Item($Row, 11).Value = $HomeOpeningDataHora
    $ws.Cells.Item($Row, 12).Value = $HomeClosingOdds
    $ws.Cells.Item($Row, 13).Value = $HomeClosingDataHora
    $ws.Cells.Item($Row, 14).Value = $DrawOpeningOdds
    $ws.Cells.Item($Row, 15).Value = $DrawOpeningDataHora
    $ws.Cells.Item($Row, 16).Value = $DrawClosingOdds
    $ws.Cells.Item($Row, 17).Value = $DrawClosingDataHora
    $ws.Cells.Item($Row, 18).Value = $AwayOpeningOdds
    $ws.Cells.Item($Row, 19).Value = $AwayOpeningDataHora
    $ws.Cells.Item($Row, 20).Value = $AwayClosingOdds
    $ws.Cells.Item($Row, 21).Value = $AwayClosingDataHora
    $ws.Cells.Item($Row, 22).Value = $Url
}

Set-MatchRow 26 25 "gibraltar" "national-league" "2023-2024" `
    45234.6875 "Glacis United" 1 "College 1975 FC" 0 `
    1.75 "04/11/2023 13:16" `
    2.1 "04/11/2023 16:08" `
    3.98 "04/11/2023 13:16" `
    3.92 "04/11/2023 16:07" `
    3.36 "04/11/2023 13:16" `
    2.67 "04/11/2023 16:08" `
    "https://www.betexplorer.com/football/gibraltar/national-league/glacis-united-college-1975/Ya5KCWUK/"

Set-MatchRow 27 26 "gibraltar" "national-league" "2023-2024" `
    45234.8125 "Mons Calpe" 0 "Lincoln Red Imps" 2 `
    9.76 "03/11/2023 19:32" `
    11.2 "04/11/2023 19:15" `
    8.73 "03/11/2023 19:32" `
    7.74 "04/11/2023 19:15" `
    1.12 "03/11/2023 19:32" `
    1.13 "04/11/2023 19:15" `
    "https://www.betexplorer.com/football/gibraltar/national-league/mons-calpe-lincoln-red-imps/Cj4OBjpR/"
